# The "2024" sheet tracks monthly transaction log entries.
# A new "Loan" entry (amazeloan, dated 2024-09-01 18:42:55) was recorded,
# which inserts a new row at the top of the "Loan" group's entries (row 10),
# shifting all subsequent rows (10-34) down by one (to 11-35).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new row at row 10, pushing existing rows 10:34 down to 11:35.
$ws.Rows("10:10").Insert()

# Populate the new row 10 with the new Loan transaction entry
# (September_Details = R, September_Date = S).
$ws.Range("R10").Value = "amazeloan"
$ws.Range("S10").Value = "2024-09-01 18:42:55"
